# Add data for 2023-11-26
# Updates violent-crime-full-year.xlsx: refresh 2023 (and a few prior-year
# correction) totals across the Citywide Totals, By Neighborhood, and each
# individual neighborhood worksheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 7281
$ws.Range("J2").Value = 6951
$ws.Range("J3").Value = 7360
$ws.Range("E4").Value = 2016
$ws.Range("J4").Value = 1599
$ws.Range("J6").Value = 9890
$ws.Range("E7").Value = 26022
$ws.Range("I7").Value = 26233
$ws.Range("J7").Value = 26376

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 120
$ws.Range("J6").Value = 199
$ws.Range("J8").Value = 1656
$ws.Range("J10").Value = 196
$ws.Range("J11").Value = 465
$ws.Range("J14").Value = 142
$ws.Range("J15").Value = 316
$ws.Range("J16").Value = 104
$ws.Range("J19").Value = 765
$ws.Range("J20").Value = 557
$ws.Range("J22").Value = 61
$ws.Range("J23").Value = 242
$ws.Range("J25").Value = 130
$ws.Range("J26").Value = 53
$ws.Range("J27").Value = 158
$ws.Range("J29").Value = 1418
$ws.Range("J31").Value = 269
$ws.Range("J33").Value = 1190
$ws.Range("J34").Value = 119
$ws.Range("J36").Value = 360
$ws.Range("J37").Value = 819
$ws.Range("J41").Value = 189
$ws.Range("J42").Value = 1142
$ws.Range("J47").Value = 194
$ws.Range("J48").Value = 299
$ws.Range("J50").Value = 158
$ws.Range("J52").Value = 676
$ws.Range("J53").Value = 387
$ws.Range("J57").Value = 120
$ws.Range("I63").Value = 184
$ws.Range("J63").Value = 84
$ws.Range("J64").Value = 174
$ws.Range("J65").Value = 662
$ws.Range("E67").Value = 1130
$ws.Range("J67").Value = 980
$ws.Range("J72").Value = 104
$ws.Range("J73").Value = 254
$ws.Range("J76").Value = 377
$ws.Range("J79").Value = 737
$ws.Range("J80").Value = 44
$ws.Range("J83").Value = 528
$ws.Range("J84").Value = 220
$ws.Range("J85").Value = 1085
$ws.Range("J86").Value = 167
$ws.Range("J88").Value = 282
$ws.Range("J94").Value = 286
$ws.Range("J96").Value = 284
$ws.Range("J97").Value = 239
$ws.Range("J100").Value = 48
$ws.Range("E101").Value = 26022
$ws.Range("I101").Value = 26233
$ws.Range("J101").Value = 26376

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 48
$ws.Range("J7").Value = 142

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 85
$ws.Range("J7").Value = 284

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 131
$ws.Range("J7").Value = 465

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 290
$ws.Range("J3").Value = 389
$ws.Range("J6").Value = 311
$ws.Range("J7").Value = 1085

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 194
$ws.Range("J6").Value = 290
$ws.Range("J7").Value = 676

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 258
$ws.Range("J7").Value = 387

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 442
$ws.Range("J6").Value = 591
$ws.Range("J7").Value = 1656

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 198
$ws.Range("J7").Value = 528

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J4").Value = 54
$ws.Range("J6").Value = 417
$ws.Range("J7").Value = 1190

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 243
$ws.Range("J6").Value = 238
$ws.Range("J7").Value = 819

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J4").Value = 28
$ws.Range("J6").Value = 243
$ws.Range("J7").Value = 662

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 65
$ws.Range("J7").Value = 269

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J3").Value = 365
$ws.Range("E4").Value = 57
$ws.Range("J4").Value = 66
$ws.Range("J6").Value = 273
$ws.Range("E7").Value = 1130
$ws.Range("J7").Value = 980

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 70
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 503
$ws.Range("J4").Value = 76
$ws.Range("J6").Value = 361
$ws.Range("J7").Value = 1418

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 47
$ws.Range("J6").Value = 146
$ws.Range("J7").Value = 299

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 185
$ws.Range("J3").Value = 219
$ws.Range("J4").Value = 36
$ws.Range("J7").Value = 765

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 377

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 60
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 228
$ws.Range("J4").Value = 47
$ws.Range("J6").Value = 607
$ws.Range("J7").Value = 1142

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 46
$ws.Range("J6").Value = 111
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 242

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 205
$ws.Range("J3").Value = 248
$ws.Range("J7").Value = 737

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 174

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 156
$ws.Range("J3").Value = 188
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 557

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 116
$ws.Range("J7").Value = 360

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 32
$ws.Range("J7").Value = 119

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J4").Value = 25
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 286

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 130

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 44
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 316

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J2").Value = 41
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J3").Value = 3
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 63
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 43
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 145
$ws.Range("J7").Value = 282

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 158

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 61

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 104

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 104
